$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 "Experimental" value: was empty -> "false" (literal text, not boolean).
# A bare Value="false" gets auto-coerced to a Boolean by Excel, and pre-setting
# NumberFormat="@" to force text creates a new style. Using a leading
# apostrophe forces literal text, then re-pasting the original cell format
# (copy/paste-special-formats from an already "style 2" cell) restores the
# untouched style index so only the value itself changes.
$ws.Range("B7").Value = "'false"
$ws.Range("B13").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# B8 "Date" value: plain text replacement, no type-coercion pitfalls.
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# B15 "Case Sensitive" value: was empty -> "true" (literal text). Same
# boolean-coercion workaround as B7.
$ws.Range("B15").Value = "'true"
$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$excel.CutCopyMode = $false
